# DataBase.xlsx update:
#  - drop the "anthony padilla" row
#  - rename "mr beast" -> "mr beast shorts"
#  - add a new column D "link to user instagram" with each channel's IG handle
#
# Operations are ordered to mirror how this was actually produced in Excel
# (rename the label first, then remove the obsolete row, then append the
# new column) so that the shared-string table / relationship ids come out
# the same way a real editing session would leave them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "mr beast" channel label before anything else moves.
$ws.Range("B4").Value = "mr beast shorts"

# 2. Hyperlinks don't follow a row delete/shift on their own here, so drop
#    them first and rebuild them afterwards against the final layout.
$ws.Hyperlinks.Delete()

# 3. Remove the "anthony padilla" row entirely (old row 3).
$ws.Rows("3:3").Delete()

# 4. Row heights shifted along with the delete - put rows 2/3 back to the
#    18pt / 14.25pt they had before (heights stay pinned to the row number).
$ws.Rows(2).RowHeight = 18
$ws.Rows(3).RowHeight = 14.25

# 5. Renumber the "sno" column (1..7).
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("A$r").Value = $r - 1
}

# 6. Re-create the hyperlinks on column C for the (now shifted) rows.
$links = @(
    "https://www.youtube.com/c/UCbAZH3nTxzyNmehmTUhuUsA",
    "https://www.youtube.com/channel/UC4-79UOlP48-QNGgCko5p2g",
    "https://www.youtube.com/c/UCA0mlN90EHCizvo101nbr-g",
    "https://www.youtube.com/channel/UCMiY4t431lhXY4QtPZtzftQ",
    "https://www.youtube.com/c/UCE9ZKI1b_PhVm3gejYuilhw",
    "https://www.youtube.com/c/UCZiJzk4wTIzaqHI4FXZ_eRQ",
    "https://www.youtube.com/c/UCS_NmOvbqaC9ccWSymx5Gpg"
)
$row = 2
foreach ($link in $links) {
    $cell = $ws.Range("C$row")
    $ws.Hyperlinks.Add($cell, $link)
    $cell.Style = "Hyperlink"
    $row = $row + 1
}

# 7. Add the new "link to user instagram" column.
$ws.Range("D1").Value = "link to user instagram"
$handles = @("sidemen", "mrbeast", "nile.red", "ksi", "impaulsiveshow", "wroetoshaw", "smosh")
$row = 2
foreach ($handle in $handles) {
    $ws.Range("D$row").Value = $handle
    $row = $row + 1
}

# Match the saved selection (active cell D1).
$ws.Range("D1").Select()
